# DO NOT TRUST PREV RESULTS! should be: 60 % err on face areas
# Updates the "missed" sheet's column A (rows 2-38) with the corrected
# face-area labels, per the commit's data re-shuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("missed")

$values = @{
    2  = "50-1-2"
    3  = "38-4-1"
    4  = "49-1-1"
    5  = "37-5-1"
    6  = "50-4-1"
    7  = "46-3-1"
    8  = "41-3-1"
    9  = "42-3-2"
    10 = "46-1-2"
    11 = "45-5-1"
    12 = "49-3-2"
    13 = "38-1-1"
    14 = "46-3-2"
    15 = "49-2-1"
    16 = "56-2-2"
    17 = "50-2-1"
    18 = "47-1-1"
    19 = "66-3-1"
    20 = "47-4-2"
    21 = "40-4-2"
    22 = "38-2-1"
    23 = "50-3-1"
    24 = "52-2-1"
    25 = "51-1-1"
    26 = "49-2-2"
    27 = "37-1-1"
    28 = "40-5-1"
    29 = "50-3-2"
    30 = "40-2-2"
    31 = "45-5-2"
    32 = "30-4-3"
    33 = "52-2-2"
    34 = "51-2-2"
    35 = "49-1-2"
    36 = "52-3-1"
    37 = "46-4-1"
    38 = "48-4-1"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
